$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ID = 1, JNE-SLD001)
$ws.Range("D2").Value = "2022-04-10 17:01:00"
$ws.Range("E2").Value = "2022-04-10 17:01:00"
$ws.Range("I2").Value = "2022-04-04 17:01:00"
$ws.Range("K2").Value = "Berlangsung"

# Row 5 (ID = 9, JNE-902109)
$ws.Range("D5").Value = "2022-04-11 16:47:00"
$ws.Range("E5").Value = "2022-04-11 16:47:00"
$ws.Range("I5").Value = "2022-04-11 16:47:00"
